# The 2007年 data row (row 2) is removed from the sheet; every row below it
# (2010年, 2012年, 2015年, 2017年) shifts up by one, which also shrinks the
# used range from A1:Y6 down to A1:Y5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
